$wb = $excel.ActiveWorkbook
$ws5 = $wb.Worksheets.Item("createAccountFormAllDataRequired")
$ws5.Activate()
$ws5.Range("C1").Value = "Password"
$ws5.Range("E3").Value = "Address "

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Sheet6"

$newSheet.Range("A1").Value = $ws5.Range("A2").Value()
$newSheet.Range("B1").Value = $ws5.Range("B2").Value()
$newSheet.Range("C1").Value = $ws5.Range("C2").Value()
$newSheet.Range("D1").Value = $ws5.Range("D2").Value()
$newSheet.Range("H1").Value = $ws5.Range("H2").Value()
$newSheet.Range("J1").Value = $ws5.Range("J2").Value()

$newSheet.Range("A2").Value = $ws5.Range("A3").Value()
$newSheet.Range("B2").Value = $ws5.Range("B3").Value()
$newSheet.Range("C2").Value = $ws5.Range("C3").Value()
$newSheet.Range("D2").Value = $ws5.Range("D3").Value()
$newSheet.Range("E2").Value = $ws5.Range("E3").Value()
$newSheet.Range("G2").Value = $ws5.Range("G3").Value()
$newSheet.Range("H2").Value = $ws5.Range("H3").Value()
$newSheet.Range("K2").Value = $ws5.Range("K3").Value()
